$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monsters")

# Rename the raid boss monster
$ws.Range("B7").Value = "Enraged and Corrupted Little Girl"

# Update stat values on row 7 (raid boss special attack stats)
$ws.Range("AI7").Value = 0.9
$ws.Range("AJ7").Value = 0.8
$ws.Range("AK7").Value = 0.9
$ws.Range("AL7").Value = 0.8
$ws.Range("AR7").Value = 0.4
$ws.Range("AV7").Value = 4
$ws.Range("AW7").Value = 0.45

# Update the active selection to match the final view state
$ws.Range("AV7").Select()
